# Sync attendance_reports: reorder "Recorded By" (column G) entries so that
# the first author in each comma-separated list is moved to the end
# (left-rotate by one), leaving single-author cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }

    $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
    $cell.Value = $rotated
}
